# "Updated RAD Scripts to remove Estate Tax and executed in Production"
#
# The RAD (Regression/Automation) test sheet records the outcome of a
# Katalon test run. The Estate Tax scenarios (rows 6-7) were pulled out of
# the active run (Result -> Fail, Execute -> DoNotRun) while the remaining
# Personal Income Tax scenarios (rows 2-5) were re-executed in Production,
# so their Date/Result timestamps were refreshed. The active selection was
# also left parked on the now-excluded Estate Tax rows (C6:C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-5 (Personal Income Tax scenarios): re-run in Production, new
# timestamps recorded. Row 2's PaymentType cell is also rewritten (same
# text as before - "Existing Liability with Notice/Invoice Number" - but
# re-pointed at the refreshed shared-string table).
$ws.Range("B2").Value = "Thu Nov 07 18:01:24 EST 2024"
$ws.Range("D2").Value = "Existing Liability with Notice/Invoice Number"

$ws.Range("B3").Value = "Thu Nov 07 18:01:37 EST 2024"

$ws.Range("B4").Value = "Thu Nov 07 18:01:51 EST 2024"

$ws.Range("B5").Value = "Thu Nov 07 18:02:06 EST 2024"

# Rows 6-7 (Estate Tax scenarios): removed from the run - marked as a
# failed/skipped execution instead of Pass/Y.
$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Thu Nov 07 16:45:28 EST 2024"
$ws.Range("C6").Value = "DoNotRun"

$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Thu Nov 07 16:45:42 EST 2024"
$ws.Range("C7").Value = "DoNotRun"

# Leave the selection on the now-excluded Estate Tax rows.
$ws.Range("C6:C7").Select()
